$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 148. This shifts all existing
# rows 148-227 down to 149-228 (Excel automatically extends the used
# range / dimension and copies formatting for the inserted row, exactly
# mirroring the target diff where every row from 148 onward contains the
# data that used to be one row above it, and a new final row 228 is
# added carrying what used to be row 227's data).
$ws.Rows.Item(148).Insert()

# Populate the newly inserted row 148 with the new weekly record.
$ws.Cells.Item(148, 1).Value = 5
$ws.Cells.Item(148, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(148, 3).Value = "Maule"
$ws.Cells.Item(148, 4).Value = 44518
$ws.Cells.Item(148, 5).Value = 7
$ws.Cells.Item(148, 6).Value = 100112023
$ws.Cells.Item(148, 7).Value = "Brócoli"
$ws.Cells.Item(148, 8).Value = "Sin especificar"
$ws.Cells.Item(148, 9).Value = "Primera"
$ws.Cells.Item(148, 10).Value = 5000
$ws.Cells.Item(148, 11).Value = 500
$ws.Cells.Item(148, 12).Value = 500
$ws.Cells.Item(148, 13).Value = 500
$ws.Cells.Item(148, 14).Value = "`$/unidad"
$ws.Cells.Item(148, 15).Value = "Región del Maule"
$ws.Cells.Item(148, 16).Value = 500
$ws.Cells.Item(148, 17).Value = 1
$ws.Cells.Item(148, 18).Value = "Hortaliza"
